# Update the departure/return date test values and refresh the
# selected cell on each sheet (API validation test data refresh).

$wb = $excel.ActiveWorkbook

$wsOneWay = $wb.Worksheets.Item("One-way")
$wsRoundTrip = $wb.Worksheets.Item("Round-trip")

# "Depart On" value moves from 22/10/2018 -> 22/12/2018 on both sheets.
$wsOneWay.Range("B4").Value = "22/12/2018"
$wsRoundTrip.Range("B4").Value = "22/12/2018"

# "Return On" value moves from 30/10/2018 -> 30/12/2018 on the round-trip sheet.
$wsRoundTrip.Range("B5").Value = "30/12/2018"

# Move the active selection on each sheet.
$wsOneWay.Activate()
$wsOneWay.Range("B4").Select()

$wsRoundTrip.Activate()
$wsRoundTrip.Range("B5").Select()
